$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated price/volume cells keep their original text representation
# (these columns store formatted strings, not numeric values), so force the
# number format to Text before writing the new values.
$cellNames = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "E15", "D16", "E16", "D17", "E17", "E18", "D20", "E20", "E21", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "D27", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "D48", "E48", "E49", "E50", "E51")
foreach ($name in $cellNames) {
    $ws.Range($name).NumberFormat = "@"
}

$updates = @{
    "D2" = "332.57"
    "E2" = "1.25%"
    "D3" = "44.16"
    "E3" = "6.35%"
    "D4" = "5.792"
    "E4" = "3.23%"
    "D5" = "0.08350"
    "E5" = "2.21%"
    "D6" = "8.817"
    "E6" = "1.09%"
    "D7" = "4.499"
    "E7" = "-0.47%"
    "D8" = "1.984"
    "E8" = "-2.00%"
    "D9" = "2.899"
    "E9" = "-1.45%"
    "D10" = "0.9340"
    "D11" = "0.1239"
    "E11" = "-2.63%"
    "D12" = "0.1956"
    "E12" = "0.31%"
    "D13" = "0.09652"
    "E13" = "2.85%"
    "D14" = "0.03940"
    "E14" = "3.51%"
    "E15" = "0.75%"
    "D16" = "0.001310"
    "E16" = "1.00%"
    "D17" = "0.006067"
    "E17" = "-2.36%"
    "E18" = "2.00%"
    "D20" = "8.964"
    "E20" = "8.32%"
    "E21" = "-1.59%"
    "E22" = "6.63%"
    "D23" = "0.04413"
    "E23" = "-0.14%"
    "D24" = "0.001260"
    "E24" = "0.08%"
    "D25" = "0.004385"
    "E25" = "0.27%"
    "E26" = "0.86%"
    "D27" = "0.0003992"
    "D39" = "0.02802"
    "E39" = "1.38%"
    "D40" = "0.05725"
    "E40" = "5.76%"
    "D41" = "0.007920"
    "E41" = "3.30%"
    "E42" = "0.98%"
    "D43" = "0.008984"
    "E43" = "0.06%"
    "D44" = "0.002102"
    "E44" = "-0.93%"
    "D45" = "0.01018"
    "E45" = "-12.04%"
    "D46" = "0.00007211"
    "E46" = "9.38%"
    "E47" = "0.01%"
    "D48" = "0.003258"
    "E48" = "0.50%"
    "E49" = "-0.05%"
    "E50" = "0.01%"
    "E51" = "0.01%"
}

foreach ($name in $updates.Keys) {
    $ws.Range($name).Value = $updates[$name]
}
